$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the commission formula (column K) for every data row (2..135) with the
# updated category rates/flat amounts:
#   RAM          65% -> 52%
#   BATTERY      70% -> 45%
#   KEYBOARD     60% -> flat 150000
#   ADAPTOR      45% -> 40%
#   MAINBOARD    15% -> flat 100000
#   ACCESSORIES  25% -> 35%
#   VGA          15% -> flat 100000
#   POWERSUPPLY  40% -> 30%
#   PROCESSOR    20% -> flat 150000
#   CASING       15% -> flat 100000
#   MONITOR       5% -> flat 100000
# (INTERNALSTORAGE, SERVICE, DISPLAY, SOFTWARE, SECOND, EXTERNALSTORAGE unchanged)
for ($r = 2; $r -le 135; $r++) {
    $formula = '=IF(E' + $r + '="INTERNALSTORAGE",H' + $r + '*52%,' + `
        'IF(E' + $r + '="SERVICE",H' + $r + '*100%,' + `
        'IF(E' + $r + '="DISPLAY",H' + $r + '*50%,' + `
        'IF(E' + $r + '="RAM",H' + $r + '*52%,' + `
        'IF(E' + $r + '="SOFTWARE",H' + $r + '*100%,' + `
        'IF(E' + $r + '="BATTERY",H' + $r + '*45%,' + `
        'IF(E' + $r + '="KEYBOARD",150000,' + `
        'IF(E' + $r + '="ADAPTOR",H' + $r + '*40%,' + `
        'IF(E' + $r + '="MAINBOARD",100000,' + `
        'IF(E' + $r + '="ACCESSORIES",H' + $r + '*35%,' + `
        'IF(E' + $r + '="VGA",100000,' + `
        'IF(E' + $r + '="POWERSUPPLY",H' + $r + '*30%,' + `
        'IF(E' + $r + '="PROCESSOR",150000,' + `
        'IF(E' + $r + '="SECOND",H' + $r + '*150%,' + `
        'IF(E' + $r + '="CASING",100000,' + `
        'IF(E' + $r + '="MONITOR",100000,' + `
        'IF(E' + $r + '="EXTERNALSTORAGE",100000,H' + $r + ')))))))))))))))))'
    $ws.Range("K" + $r).Formula = $formula
}

# New accounting-style number format (same as the existing IDR format but without
# the "??" alignment placeholders before the closing paren) applied to the whole
# commission column range in one shot so every cell shares a single new style.
$ws.Range("K2:K135").NumberFormat = '_([$IDR]\ * #,##0_);_([$IDR]\ * \(#,##0\);_([$IDR]\ * "-"_);_(@_)'

# Size column K to fit its new contents and select it the way the author left it.
$ws.Columns("K").ColumnWidth = 11.6
$ws.Range("K2:K135").Select()
